$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(259752.536700512, 281922.012947309, 274248.378060604, 245256.69534042,  237583.060453715),
    @(251989.722981694, 277370.710163781, 268585.458472162, 235393.987491225, 226608.735799606),
    @(292942.357805024, 326352.862188839, 314788.312392175, 271096.403217873, 259531.853421208),
    @(277055.486191765, 311960.756977041, 299878.815994494, 254232.156389037, 242150.21540649),
    @(315137.539827499, 358277.584917803, 343345.299132892, 286929.780522106, 271997.494737194),
    @(320711.544772731, 367855.886465407, 351537.572893343, 289885.516652119, 273567.203080055),
    @(321350.759292527, 371629.614500432, 354226.335636651, 288475.182948402, 271071.904084621),
    @(282003.197635088, 328643.360015535, 312499.560623076, 251506.834647099, 235363.03525464),
    @(245135.56044965,  287756.019509858, 273003.580711082, 217267.540188218, 202515.101389442),
    @(239343.691351815, 282892.845437558, 267818.952650978, 210868.430052652, 195794.537266072),
    @(251146.856964712, 298789.678596372, 282298.823598091, 219994.890331332, 203504.035333051),
    @(256375.940191033, 306920.66738509,  289425.360882694, 223326.519499371, 205831.212996975)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
